# Invalid Login page script + updated Login page data
# Rename Sheet1 -> ValidLogin, populate valid credentials,
# add a new InvalidLogin sheet with invalid credentials.

$wb = $excel.ActiveWorkbook

# --- ValidLogin (existing sheet renamed & filled) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ValidLogin"
$ws1.Range("A1").Value = "UserName"
$ws1.Range("B1").Value = "Password"
$ws1.Range("A2").Value = "admin"
$ws1.Range("B2").Value = "manager"

# --- InvalidLogin (new sheet, placed right after ValidLogin) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "InvalidLogin"
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "Bhanu"
$ws2.Range("B2").Value = "Damager"

# Auto-fit the columns on InvalidLogin to match their content widths
$null = $ws2.Columns.Item(1).AutoFit()
$null = $ws2.Columns.Item(2).AutoFit()

# Match view state: selection + zoom per sheet
$null = $ws1.Select()
$null = $ws1.Range("A1:B2").Select()
$excel.ActiveWindow.Zoom = 235

$null = $ws2.Select()
$null = $ws2.Range("B3").Select()
$excel.ActiveWindow.Zoom = 250

# InvalidLogin ends up the active (visible) tab
$ws2.Activate()
